$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link URL) -----------------------------------
# These values are not numeric-looking, so plain .Value assignment keeps them
# stored as text without needing any number-format coercion.
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'

# --- Numeric-looking columns (Price / Volume / Hora) ------------------------
# Force text format first so values such as "327.74", "-0.29%" or "13" are
# stored as literal text (matching the source data) instead of being parsed
# into numbers/percentages by Excel.
$numericLikeCells = @{
    "D2" = '327.74'
    "E2" = '-0.29%'
    "G2" = '13'
    "D3" = '44.08'
    "E3" = '2.30%'
    "G3" = '13'
    "D4" = '5.575'
    "E4" = '-0.49%'
    "G4" = '13'
    "D5" = '0.08068'
    "E5" = '-1.70%'
    "G5" = '13'
    "D6" = '1.915'
    "E6" = '0.07%'
    "G6" = '13'
    "D7" = '4.293'
    "E7" = '-3.85%'
    "G7" = '13'
    "D8" = '0.9468'
    "E8" = '0.06%'
    "G8" = '13'
    "D9" = '2.534'
    "E9" = '-11.26%'
    "G9" = '13'
    "D10" = '0.1177'
    "E10" = '-3.52%'
    "G10" = '13'
    "D11" = '0.1844'
    "E11" = '-4.06%'
    "G11" = '13'
    "D12" = '0.09663'
    "E12" = '-1.76%'
    "G12" = '13'
    "D13" = '0.04386'
    "E13" = '-1.63%'
    "G13" = '13'
    "D14" = '0.1069'
    "E14" = '-0.15%'
    "G14" = '13'
    "E15" = '0.17%'
    "G15" = '13'
    "D16" = '0.005900'
    "E16" = '-3.76%'
    "G16" = '13'
    "D17" = '3.503'
    "E17" = '0.50%'
    "G17" = '13'
    "D18" = '0.3496'
    "E18" = '-1.16%'
    "G18" = '13'
    "D19" = '9.605'
    "E19" = '9.07%'
    "G19" = '13'
    "D20" = '0.1362'
    "E20" = '-0.09%'
    "G20" = '13'
    "D21" = '0.2650'
    "E21" = '-2.97%'
    "G21" = '13'
    "D22" = '0.04215'
    "E22" = '-4.34%'
    "G22" = '13'
    "D23" = '0.001246'
    "E23" = '0.06%'
    "G23" = '13'
    "D24" = '0.004482'
    "E24" = '1.97%'
    "G24" = '13'
    "D25" = '0.0001262'
    "E25" = '1.93%'
    "G25" = '13'
    "D26" = '0.0003995'
    "E26" = '-0.49%'
    "G26" = '13'
    "G27" = '13'
    "G28" = '13'
    "G29" = '13'
    "G30" = '13'
    "G31" = '13'
    "G32" = '13'
    "G33" = '13'
    "G34" = '13'
    "G35" = '13'
    "G36" = '13'
    "G37" = '13'
    "D38" = '0.02658'
    "E38" = '-4.16%'
    "G38" = '13'
    "D39" = '0.05509'
    "E39" = '-4.07%'
    "G39" = '13'
    "E40" = '-4.08%'
    "G40" = '13'
    "D41" = '0.1398'
    "E41" = '-1.41%'
    "G41" = '13'
    "E42" = '-28.86%'
    "G42" = '13'
    "D43" = '0.002011'
    "E43" = '-3.47%'
    "G43" = '13'
    "D44" = '0.008379'
    "E44" = '-14.02%'
    "G44" = '13'
    "D45" = '0.00006896'
    "E45" = '-5.24%'
    "G45" = '13'
    "E46" = '-0.50%'
    "G46" = '13'
    "D47" = '0.002274'
    "E47" = '-0.49%'
    "G47" = '13'
    "D48" = '0.005658'
    "E48" = '67.30%'
    "G48" = '13'
    "D49" = '0.00002103'
    "E49" = '-0.50%'
    "G49" = '13'
    "D50" = '0.0002003'
    "E50" = '-0.50%'
    "G50" = '13'
    "G51" = '13'
}
foreach ($addr in $numericLikeCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLikeCells[$addr]
}
